$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '61.139.80'
$ws.Range('E2').Value = '  -4.90%  '

Set-TextValue 'D3' '3.312.96'
$ws.Range('E3').Value = '  -5.19%  '

Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  -0.08%  '

Set-TextValue 'D5' '569.76'
$ws.Range('E5').Value = '  -3.22%  '

Set-TextValue 'D6' '126.85'
$ws.Range('E6').Value = '  -5.12%  '

$ws.Range('E7').Value = '  +0.00%  '

Set-TextValue 'D8' '3.312.73'
$ws.Range('E8').Value = '  -5.15%  '

$ws.Range('E9').Value = '  -2.02%  '

Set-TextValue 'D10' '7.30'
$ws.Range('E10').Value = '  -4.29%  '

$ws.Range('E11').Value = '  -5.87%  '

Set-TextValue 'D12' '0.373'
$ws.Range('E12').Value = '  -3.75%  '

Set-TextValue 'D13' '3.877.33'
$ws.Range('E13').Value = '  -5.19%  '

$ws.Range('E14').Value = '  -1.43%  '

Set-TextValue 'D15' '3.313.81'
$ws.Range('E15').Value = '  -5.20%  '

$ws.Range('E16').Value = '  -6.97%  '

Set-TextValue 'D17' '61.210.67'
$ws.Range('E17').Value = '  -4.69%  '

Set-TextValue 'D18' '24.46'
$ws.Range('E18').Value = '  -3.31%  '

Set-TextValue 'D19' '5.59'
$ws.Range('E19').Value = '  -3.25%  '

Set-TextValue 'D20' '9.04'
$ws.Range('E20').Value = '  -10.10%  '

Set-TextValue 'D21' '13.15'
$ws.Range('E21').Value = '  -2.98%  '

Set-TextValue 'D22' '351.36'
$ws.Range('E22').Value = '  -8.99%  '

$ws.Range('E23').Value = '  -5.07%  '

$ws.Range('E24').Value = '  +0.03%  '

Set-TextValue 'D25' '3.445.92'
$ws.Range('E25').Value = '  -5.17%  '

Set-TextValue 'D26' '69.80'
$ws.Range('E26').Value = '  -5.80%  '

$ws.Range('E27').Value = '  -7.38%  '

Set-TextValue 'D28' '0.998'
$ws.Range('E28').Value = '  -0.12%  '

Set-TextValue 'D29' '7.19'
$ws.Range('E29').Value = '  -2.36%  '

$ws.Range('E30').Value = '  -3.87%  '

Set-TextValue 'D31' '7.82'
$ws.Range('E31').Value = '  -4.04%  '

$ws.Range('E32').Value = '  -6.49%  '

$ws.Range('E33').Value = '  -0.03%  '

Set-TextValue 'D34' '0.148'
$ws.Range('E34').Value = '  -4.64%  '

Set-TextValue 'D35' '3.344.50'
$ws.Range('E35').Value = '  -5.08%  '

Set-TextValue 'D36' '22.36'
$ws.Range('E36').Value = '  -4.05%  '

Set-TextValue 'D37' '5.30'
$ws.Range('E37').Value = '  -0.46%  '

$ws.Range('E38').Value = '  -3.91%  '

Set-TextValue 'D39' '162.29'
$ws.Range('E39').Value = '  -2.00%  '

$ws.Range('E40').Value = '  -4.07%  '

Set-TextValue 'D41' '0.0749'
$ws.Range('E41').Value = '  -4.49%  '

Set-TextValue 'D42' '1.00'
$ws.Range('E42').Value = '  +0.02%  '

Set-TextValue 'D43' '41.01'
$ws.Range('E43').Value = '  -2.06%  '

Set-TextValue 'D44' '0.743'
$ws.Range('E44').Value = '  -7.87%  '

$ws.Range('E45').Value = '  -4.58%  '

$ws.Range('E46').Value = '  -5.55%  '

$ws.Range('E47').Value = '  -6.07%  '

Set-TextValue 'D48' '22.27'
$ws.Range('E48').Value = '  -8.93%  '

Set-TextValue 'D49' '6.63'
$ws.Range('E49').Value = '  -2.63%  '

Set-TextValue 'D50' '0.854'
$ws.Range('E50').Value = '  -6.50%  '

Set-TextValue 'D51' '2.199.62'
$ws.Range('E51').Value = '  -9.70%  '
